$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header value for new column C (t+3)
$ws.Range("C1").Value = 2

# Copy the style of B1 (bordered/bold header style) onto C1
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Data values for column C, rows 2-10
$values = @(
    -5.030479892299043,
    -1.153235334261761,
    -0.07292569980107731,
    -0.417749988516372,
    0.01457436480836208,
    0.1039564587721915,
    0.1352692197136115,
    0.02838054686518928,
    0.02323175602225529
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
